$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "Type"
$ws.Range("H1").Value = "Sources"

# Seed the shared-string table with the new distinct values in the same
# order they were introduced by the original author (Type options then
# Source options, interleaved), using a scratch area that is cleared
# again right after. This keeps the workbook's string table ordering
# consistent with the source edit.
$scratch = $ws.Range("J1:J7")
$scratch.Item(1).Value = "Reporting"
$scratch.Item(2).Value = "Reposting"
$scratch.Item(3).Value = "Activist"
$scratch.Item(4).Value = "Internal"
$scratch.Item(5).Value = "Opinion"
$scratch.Item(6).Value = "Crowdsourced"
$scratch.Item(7).Value = "Other"
$scratch.ClearContents()

# Data rows: Name (for reference) -> Type (G), Sources (H)
$data = @(
    @{ Row = 2;  G = "Activist";  H = "Internal" },
    @{ Row = 3;  G = "Reporting"; H = "Internal" },
    @{ Row = 4;  G = "Activist";  H = "Crowdsourced" },
    @{ Row = 5;  G = "Activist";  H = "Internal" },
    @{ Row = 6;  G = "Activist";  H = "Crowdsourced" },
    @{ Row = 7;  G = "Reporting"; H = "Internal" },
    @{ Row = 8;  G = "Reporting"; H = "Internal" },
    @{ Row = 9;  G = "Reporting"; H = "Internal" },
    @{ Row = 10; G = "Opinion";   H = "Reposting" },
    @{ Row = 11; G = "Opinion";   H = "Reposting" },
    @{ Row = 12; G = "Reporting"; H = "Reposting" },
    @{ Row = 13; G = "Reporting"; H = "Internal" },
    @{ Row = 14; G = "Reporting"; H = "Internal" },
    @{ Row = 15; G = "Activist";  H = "Crowdsourced" },
    @{ Row = 17; G = "Opinion";   H = "Internal" },
    @{ Row = 18; G = "Other";     H = "Other" },
    @{ Row = 19; G = "Other";     H = "Other" },
    @{ Row = 20; G = "Other";     H = "Other" },
    @{ Row = 21; G = "Reporting"; H = "Crowdsourced" },
    @{ Row = 22; G = "Opinion";   H = "Reposting" },
    @{ Row = 23; G = "Activist";  H = "Reposting" },
    @{ Row = 24; G = "Opinion";   H = "Internal" },
    @{ Row = 25; G = "Other";     H = "Internal" },
    @{ Row = 26; G = "Reporting"; H = "Crowdsourced" },
    @{ Row = 27; G = "Reporting"; H = "Reposting" },
    @{ Row = 28; G = "Reporting"; H = "Reposting" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}

# Column widths for the new columns (bestFit-equivalent, matching the
# widths Excel computed for the "language"/"Type"/"Sources" columns).
$ws.Columns("F").ColumnWidth = 9.529947916666666
$ws.Columns("G").ColumnWidth = 8.072916666666666
$ws.Columns("H").ColumnWidth = 8.166666666666666

# Update the view: move the selection to match the author's final position
$ws.Range("G21").Select() | Out-Null
